# cambios de agosto, puntos fe de ratas e historico
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 data updates (Q2 2022 report replaces Q1 2022 report) ---

# Reporting period start/end dates
$ws.Range("B8").Value2 = 44652   # 2022-04-01
$ws.Range("C8").Value2 = 44742   # 2022-06-30

# Denominación de cada informe / Área responsable
$ws.Range("D8").Value2 = "Primera Sesión Ordinaria 2022"
$ws.Range("E8").Value2 = "Abogado General"

# Fecha en que se presentó y/o entregó el informe
$ws.Range("H8").Value2 = 44693   # 2022-05-12

# Hipervínculo al documento del informe correspondiente
$ws.Range("I8").Value2 = "http://transparenciadocs.hidalgo.gob.mx/ENTIDADES/UPPachuca/dir1/2022/Abril-Junio/29/Informe%20de%20actividades%20del%20Titular.pdf"

# Fecha de validación / Fecha de actualización
$ws.Range("K8").Value2 = 44753   # 2022-07-11
$ws.Range("L8").Value2 = 44753   # 2022-07-11

# --- Sheet view / selection ---
[void]$ws.Range("H13").Select()
